$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '60.951.61'
Set-TextValue $ws.Range("E2") '  +0.00%  '

Set-TextValue $ws.Range("D3") '2.920.12'

Set-TextValue $ws.Range("E4") '  -0.04%  '

Set-TextValue $ws.Range("D5") '590.35'
Set-TextValue $ws.Range("E5") '  +0.66%  '

Set-TextValue $ws.Range("D6") '146.62'
Set-TextValue $ws.Range("E6") '  +0.81%  '

Set-TextValue $ws.Range("E7") '  -0.04%  '

Set-TextValue $ws.Range("E8") '  +0.15%  '

Set-TextValue $ws.Range("D9") '6.90'
Set-TextValue $ws.Range("E9") '  -0.05%  '

Set-TextValue $ws.Range("E10") '  -0.72%  '

Set-TextValue $ws.Range("E11") '  -1.45%  '

Set-TextValue $ws.Range("E12") '  +0.07%  '

Set-TextValue $ws.Range("D13") '33.59'
Set-TextValue $ws.Range("E13") '  +0.01%  '

Set-TextValue $ws.Range("E14") '  +0.02%  '

Set-TextValue $ws.Range("D15") '3.403.37'
Set-TextValue $ws.Range("E15") '  +0.00%  '

Set-TextValue $ws.Range("D16") '60.898.43'
Set-TextValue $ws.Range("E16") '  +0.00%  '

Set-TextValue $ws.Range("E17") '  -0.93%  '

Set-TextValue $ws.Range("D18") '2.919.88'
Set-TextValue $ws.Range("E18") '  -0.04%  '

Set-TextValue $ws.Range("D19") '432.65'
Set-TextValue $ws.Range("E19") '  +0.73%  '

Set-TextValue $ws.Range("D20") '13.41'
Set-TextValue $ws.Range("E20") '  -1.47%  '

Set-TextValue $ws.Range("D21") '0.678'
Set-TextValue $ws.Range("E21") '  -0.55%  '

Set-TextValue $ws.Range("D22") '7.12'
Set-TextValue $ws.Range("E22") '  -0.22%  '

Set-TextValue $ws.Range("D23") '81.28'
Set-TextValue $ws.Range("E23") '  +0.97%  '

Set-TextValue $ws.Range("D24") '10.89'
Set-TextValue $ws.Range("E24") '  +1.29%  '

Set-TextValue $ws.Range("E25") '  -0.94%  '

Set-TextValue $ws.Range("D26") '11.86'
Set-TextValue $ws.Range("E26") '  -0.87%  '

Set-TextValue $ws.Range("E27") '  -0.02%  '

Set-TextValue $ws.Range("D28") '2.30'
Set-TextValue $ws.Range("E28") '  +6.18%  '

Set-TextValue $ws.Range("E29") '  -0.28%  '

Set-TextValue $ws.Range("D30") '6.98'
Set-TextValue $ws.Range("E30") '  -2.85%  '

Set-TextValue $ws.Range("E31") '  +3.35%  '

Set-TextValue $ws.Range("D32") '26.65'
Set-TextValue $ws.Range("E32") '  +0.29%  '

Set-TextValue $ws.Range("E33") '  +0.00%  '

Set-TextValue $ws.Range("D34") '0.0₃0865'
Set-TextValue $ws.Range("E34") '  -0.40%  '

Set-TextValue $ws.Range("E35") '  -0.10%  '

Set-TextValue $ws.Range("E36") '  -0.16%  '

Set-TextValue $ws.Range("D37") '3.01'
Set-TextValue $ws.Range("E37") '  -0.17%  '

Set-TextValue $ws.Range("E38") '  -1.24%  '

Set-TextValue $ws.Range("E39") '  -4.80%  '

Set-TextValue $ws.Range("D40") '8.55'
Set-TextValue $ws.Range("E40") '  -0.96%  '

Set-TextValue $ws.Range("D41") '41.42'
Set-TextValue $ws.Range("E41") '  +0.31%  '

Set-TextValue $ws.Range("D42") '0.283'
Set-TextValue $ws.Range("E42") '  -4.25%  '

Set-TextValue $ws.Range("D43") '376.50'
Set-TextValue $ws.Range("E43") '  -0.50%  '

Set-TextValue $ws.Range("D44") '2.705.81'
Set-TextValue $ws.Range("E44") '  +0.21%  '

Set-TextValue $ws.Range("E45") '  -2.04%  '

Set-TextValue $ws.Range("D46") '133.90'
Set-TextValue $ws.Range("E46") '  +0.73%  '

Set-TextValue $ws.Range("E47") '  +0.04%  '

Set-TextValue $ws.Range("D48") '23.93'
Set-TextValue $ws.Range("E48") '  -3.43%  '

Set-TextValue $ws.Range("E49") '  -0.50%  '

Set-TextValue $ws.Range("E50") '  -2.67%  '

Set-TextValue $ws.Range("E51") '  -0.63%  '
